$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: replace Firstname/Lastname (Harley -> Michael, Dion -> Watt)
$ws.Range("B2").Value = "Michael"
$ws.Range("C2").Value = "Watt"

# Row 3: replace Firstname (Xinjon -> Bazun)
$ws.Range("B3").Value = "Bazun"

# Row 5 (new record)
$ws.Range("A5").Value = "Miss"
$ws.Range("B5").Value = "Zia"
$ws.Range("C5").Value = "Gill"
$ws.Range("D5").Value = "Apple"

# Row 6 (new record)
$ws.Range("A6").Value = "Sir"
$ws.Range("B6").Value = "Sean"
$ws.Range("C6").Value = "Paul"
$ws.Range("D6").Value = "Persistent"

# Update selection to match the new active range
$ws.Range("A2:D6").Select()
